{"js": "const pairs = [\n  [\"N = 89,510\", \"N = 87,276\"],\n  [\"96.4 (89.8, 102.7)\", \"96.5 (89.8, 102.7)\"],\n  [\"463.3 (319.2, 640.3)\", \"462.2 (318.4, 638.1)\"],\n  [\"233.6 (114.5, 403.0)\", \"232.5 (114.0, 401.1)\"],\n  [\"585.3 (372.1, 851.8)\", \"584.2 (371.7, 849.1)\"],\n  [\"385.0 (140.0, 630.0)\", \"350.0 (140.0, 630.0)\"],\n  [\"1,984 (2.2)\", \"1,917 (2.2)\"],\n  [\"1,274 (1.4)\", \"1,227 (1.4)\"],\n  [\"2,659 (3.0)\", \"2,590 (3.0)\"],\n  [\"86,851 (97)\", \"84,686 (97)\"],\n  [\"51,446 (57)\", \"50,371 (58)\"],\n  [\"38,064 (43)\", \"36,905 (42)\"],\n  [\"6,980 (7.8)\", \"6,683 (7.7)\"],\n  [\"21,885 (24)\", \"21,174 (24)\"],\n  [\"21,111 (24)\", \"20,530 (24)\"],\n  [\"39,534 (44)\", \"38,889 (45)\"],\n  [\"11,436 (13)\", \"11,034 (13)\"],\n  [\"19,426 (22)\", \"18,743 (21)\"],\n  [\"23,495 (26)\", \"22,949 (26)\"],\n  [\"20,810 (23)\", \"20,512 (24)\"],\n  [\"6,061 (6.8)\", \"6,000 (6.9)\"],\n  [\"8,282 (9.3)\", \"8,038 (9.2)\"],\n  [\"51,741 (58)\", \"50,562 (58)\"],\n  [\"31,776 (35)\", \"30,886 (35)\"],\n  [\"5,993 (6.7)\", \"5,828 (6.7)\"],\n  [\"4,901 (5.5)\", \"4,765 (5.5)\"],\n  [\"18,087 (20)\", \"17,618 (20)\"],\n  [\"22,516 (25)\", \"21,962 (25)\"],\n  [\"23,479 (26)\", \"22,953 (26)\"],\n  [\"20,527 (23)\", \"19,978 (23)\"],\n  [\"63,794 (71)\", \"62,223 (71)\"],\n  [\"22,569 (25)\", \"21,985 (25)\"],\n  [\"3,147 (3.5)\", \"3,068 (3.5)\"],\n  [\"16,138 (18)\", \"15,701 (18)\"],\n  [\"30,358 (34)\", \"29,625 (34)\"],\n  [\"43,014 (48)\", \"41,950 (48)\"],\n  [\"73,949 (83)\", \"72,160 (83)\"],\n  [\"14,854 (17)\", \"14,426 (17)\"],\n  [\"707 (0.8)\", \"690 (0.8)\"],\n  [\"75,936 (85)\", \"74,100 (85)\"],\n  [\"13,156 (15)\", \"12,770 (15)\"],\n  [\"418 (0.5)\", \"406 (0.5)\"],\n  [\"29,195 (33)\", \"28,516 (33)\"],\n  [\"53,648 (60)\", \"52,300 (60)\"],\n  [\"6,667 (7.4)\", \"6,460 (7.4)\"],\n  [\"15,255 (17)\", \"14,915 (17)\"],\n  [\"65,969 (74)\", \"64,378 (74)\"],\n  [\"8,286 (9.3)\", \"7,983 (9.1)\"],\n];\n\nconst body = context.document.body;\nlet notFound = [];\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load('items');\n  await context.sync();\n  if (results.items.length === 0) {\n    notFound.push(oldText);\n    continue;\n  }\n  for (const item of results.items) {\n    item.insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\nif (notFound.length > 0) {\n  throw new Error('Not found: ' + JSON.stringify(notFound));\n}\nreturn 'done';", "ps1": "$d = $word.ActiveDocument\n\n$pairs = @(\n  @('N = 89,510', 'N = 87,276'),\n  @('96.4 (89.8, 102.7)', '96.5 (89.8, 102.7)'),\n  @('463.3 (319.2, 640.3)', '462.2 (318.4, 638.1)'),\n  @('233.6 (114.5, 403.0)', '232.5 (114.0, 401.1)'),\n  @('585.3 (372.1, 851.8)', '584.2 (371.7, 849.1)'),\n  @('385.0 (140.0, 630.0)', '350.0 (140.0, 630.0)'),\n  @('1,984 (2.2)', '1,917 (2.2)'),\n  @('1,274 (1.4)', '1,227 (1.4)'),\n  @('2,659 (3.0)', '2,590 (3.0)'),\n  @('86,851 (97)', '84,686 (97)'),\n  @('51,446 (57)', '50,371 (58)'),\n  @('38,064 (43)', '36,905 (42)'),\n  @('6,980 (7.8)', '6,683 (7.7)'),\n  @('21,885 (24)', '21,174 (24)'),\n  @('21,111 (24)', '20,530 (24)'),\n  @('39,534 (44)', '38,889 (45)'),\n  @('11,436 (13)', '11,034 (13)'),\n  @('19,426 (22)', '18,743 (21)'),\n  @('23,495 (26)', '22,949 (26)'),\n  @('20,810 (23)', '20,512 (24)'),\n  @('6,061 (6.8)', '6,000 (6.9)'),\n  @('8,282 (9.3)', '8,038 (9.2)'),\n  @('51,741 (58)', '50,562 (58)'),\n  @('31,776 (35)', '30,886 (35)'),\n  @('5,993 (6.7)', '5,828 (6.7)'),\n  @('4,901 (5.5)', '4,765 (5.5)'),\n  @('18,087 (20)', '17,618 (20)'),\n  @('22,516 (25)', '21,962 (25)'),\n  @('23,479 (26)', '22,953 (26)'),\n  @('20,527 (23)', '19,978 (23)'),\n  @('63,794 (71)', '62,223 (71)'),\n  @('22,569 (25)', '21,985 (25)'),\n  @('3,147 (3.5)', '3,068 (3.5)'),\n  @('16,138 (18)', '15,701 (18)'),\n  @('30,358 (34)', '29,625 (34)'),\n  @('43,014 (48)', '41,950 (48)'),\n  @('73,949 (83)', '72,160 (83)'),\n  @('14,854 (17)', '14,426 (17)'),\n  @('707 (0.8)', '690 (0.8)'),\n  @('75,936 (85)', '74,100 (85)'),\n  @('13,156 (15)', '12,770 (15)'),\n  @('418 (0.5)', '406 (0.5)'),\n  @('29,195 (33)', '28,516 (33)'),\n  @('53,648 (60)', '52,300 (60)'),\n  @('6,667 (7.4)', '6,460 (7.4)'),\n  @('15,255 (17)', '14,915 (17)'),\n  @('65,969 (74)', '64,378 (74)'),\n  @('8,286 (9.3)', '7,983 (9.1)'),\n)\n\n$notFound = @()\nforeach ($pair in $pairs) {\n  $oldText = $pair[0]\n  $newText = $pair[1]\n  $find = $d.Content.Find\n  $find.ClearFormatting()\n  $find.Replacement.ClearFormatting()\n  $find.Text = $oldText\n  $find.Replacement.Text = $newText\n  $find.Forward = $true\n  $find.Wrap = 1\n  $found = $find.Execute($null, $true, $false, $false, $false, $false, $true, 1, $false, $null, 2)\n  if (-not $found) {\n    $notFound += $oldText\n  }\n}\n\nif ($notFound.Count -gt 0) {\n  throw \"Not found: \" + ($notFound -join \", \")\n}"}
